$wb = $excel.ActiveWorkbook

# --- Sheet 1: LH_WF_REGISTRATION_REVIEW ---
$ws1 = $wb.Worksheets.Item(1)

# Fill in row 3 with the new review entry (column order matches row 1/2:
# A=date, B=ID, C=Reviewer, D=Version, E=Review Comments, F=Actions,
# G=Owner, H=Owner Status, I=Reviewer verification)
$ws1.Range("B3").Value = "LH_WF_REGISTRATION_REVIEW_002"
$ws1.Range("A3").Value = "27/4/2025"
$ws1.Range("C3").Value = "Ahmed Abuzaid"
$ws1.Range("D3").Value = "v1.1"
$ws1.Range("E3").Value = "back to SRS I found there are many error messages for many validations like existing user name or existing email or validations fro password"
$ws1.Range("F3").Value = "so I prefer to add some error messages ""with red color"" from the SRS to wireframe to make it more expressive, you can back to login wireframe to understand what I mean"
$ws1.Range("G3").Value = "Gehad"
$ws1.Range("H3").Value = "open"
$ws1.Range("I3").Value = "open"

# Row 3 grows tall to fit the wrapped review-comment text
$ws1.Rows.Item(3).RowHeight = 150

# Move the sheet's selection/scroll to I3 (without leaving this sheet
# marked as the active tab, since VERSION-HISTORY stays the active sheet)
$ws1.Activate()
$ws1.Range("I3").Select()

# --- Sheet 2: VERSION-HISTORY ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()

$ws2.Range("A4").Value = "v1.2"
$ws2.Range("B4").Value = "Ahmed Abuzaid"
$ws2.Range("C4").Value = "ask to add more details to registration form wireframe"
$ws2.Range("D4").Value = "4/27/2025"

$ws2.Rows.Item(4).RowHeight = 37.5

$ws2.Range("C4").Select()
